$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 26.29132066666667
$ws.Range("H2").Value = 78.873962
$ws.Range("I2").Value = 0.1411782207947891
$ws.Range("J2").Value = 0.1411782207947891
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 43.91845900000001
$ws.Range("N2").Value = 131.755377
$ws.Range("O2").Value = 0.8150909120558799
$ws.Range("P2").Value = 0.81509091205588
$ws.Range("Q2").Value = 1154.674288754853
$ws.Range("R2").Value = 10392.06859879368
$ws.Range("S2").Value = 0.115073084750051
$ws.Range("T2").Value = 0.115073084750051

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 26.29132066666667
$ws.Range("H3").Value = 78.873962
$ws.Range("I3").Value = 0.1411782207947891
$ws.Range("J3").Value = 0.1411782207947891
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.100310333333333
$ws.Range("N3").Value = 6.300930999999999
$ws.Range("O3").Value = 0.03898005312975703
$ws.Range("P3").Value = 0.03898005312975703
$ws.Range("Q3").Value = 55.21993247318022
$ws.Range("R3").Value = 496.979392258622
$ws.Range("S3").Value = 0.005503134547345446
$ws.Range("T3").Value = 0.005503134547345446

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 26.29132066666667
$ws.Range("H4").Value = 78.873962
$ws.Range("I4").Value = 0.1411782207947891
$ws.Range("J4").Value = 0.1411782207947891
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.8629
$ws.Range("N4").Value = 23.5887
$ws.Range("O4").Value = 0.1459290348143631
$ws.Range("P4").Value = 0.1459290348143631
$ws.Range("Q4").Value = 206.7260252699333
$ws.Range("R4").Value = 1860.5342274294
$ws.Range("S4").Value = 0.02060200149739261
$ws.Range("T4").Value = 0.02060200149739261

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 143.4723713333333
$ws.Range("H5").Value = 430.417114
$ws.Range("I5").Value = 0.7704129577533824
$ws.Range("J5").Value = 0.7704129577533824
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 43.91845900000001
$ws.Range("N5").Value = 131.755377
$ws.Range("O5").Value = 0.8150909120558799
$ws.Range("P5").Value = 0.81509091205588
$ws.Range("Q5").Value = 6301.085458035775
$ws.Range("R5").Value = 56709.76912232198
$ws.Range("S5").Value = 0.6279566003948726
$ws.Range("T5").Value = 0.6279566003948726

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 143.4723713333333
$ws.Range("H6").Value = 430.417114
$ws.Range("I6").Value = 0.7704129577533824
$ws.Range("J6").Value = 0.7704129577533824
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.100310333333333
$ws.Range("N6").Value = 6.300930999999999
$ws.Range("O6").Value = 0.03898005312975703
$ws.Range("P6").Value = 0.03898005312975703
$ws.Range("Q6").Value = 301.3365040592371
$ws.Range("R6").Value = 2712.028536533134
$ws.Range("S6").Value = 0.0300307380250801
$ws.Range("T6").Value = 0.0300307380250801

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 143.4723713333333
$ws.Range("H7").Value = 430.417114
$ws.Range("I7").Value = 0.7704129577533824
$ws.Range("J7").Value = 0.7704129577533824
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.8629
$ws.Range("N7").Value = 23.5887
$ws.Range("O7").Value = 0.1459290348143631
$ws.Range("P7").Value = 0.1459290348143631
$ws.Range("Q7").Value = 1128.108908556866
$ws.Range("R7").Value = 10152.9801770118
$ws.Range("S7").Value = 0.1124256193334298
$ws.Range("T7").Value = 0.1124256193334298

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 16.46418733333334
$ws.Range("H8").Value = 49.39256200000001
$ws.Range("I8").Value = 0.08840882145182853
$ws.Range("J8").Value = 0.08840882145182853
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 43.91845900000001
$ws.Range("N8").Value = 131.755377
$ws.Range("O8").Value = 0.8150909120558799
$ws.Range("P8").Value = 0.81509091205588
$ws.Range("Q8").Value = 723.0817363673195
$ws.Range("R8").Value = 6507.735627305875
$ws.Range("S8").Value = 0.07206122691095636
$ws.Range("T8").Value = 0.07206122691095637

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 16.46418733333334
$ws.Range("H9").Value = 49.39256200000001
$ws.Range("I9").Value = 0.08840882145182853
$ws.Range("J9").Value = 0.08840882145182853
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.100310333333333
$ws.Range("N9").Value = 6.300930999999999
$ws.Range("O9").Value = 0.03898005312975703
$ws.Range("P9").Value = 0.03898005312975703
$ws.Range("Q9").Value = 34.57990278613578
$ws.Range("R9").Value = 311.219125075222
$ws.Range("S9").Value = 0.003446180557331479
$ws.Range("T9").Value = 0.003446180557331479

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 16.46418733333334
$ws.Range("H10").Value = 49.39256200000001
$ws.Range("I10").Value = 0.08840882145182853
$ws.Range("J10").Value = 0.08840882145182853
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.8629
$ws.Range("N10").Value = 23.5887
$ws.Range("O10").Value = 0.1459290348143631
$ws.Range("P10").Value = 0.1459290348143631
$ws.Range("Q10").Value = 129.4562585832667
$ws.Range("R10").Value = 1165.1063272494
$ws.Range("S10").Value = 0.0129014139835407
$ws.Range("T10").Value = 0.0129014139835407
